# Generate Report for Archive
# Update status text from "Ready for handoff" to "In Translation" on every
# sheet, and shrink the now-shorter text's columns to match.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Target OOXML column <col width="..."> after the edit; the COM ColumnWidth
# setter takes "characters" and is offset from the stored width by the
# standard 5/6-character cell-padding constant, so back that out here.
$targetStoredWidth = 13.4101845877511
$newColumnWidth = $targetStoredWidth - (5 / 6)

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ($oldStatus -eq $cell.Value()) {
            $cell.Value = $newStatus
            $ws.Columns.Item($cell.Column).ColumnWidth = $newColumnWidth
        }
    }
}
